$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New OHLCV rows to append (row, A=timestamp, B=open, C=high, D=low, E=close, F=volume)
$newRows = @(
    @(1169, 45534.5, 138.65, 140.73, 137.9, 139.8, 433361.888),
    @(1170, 45534.58333333334, 139.81, 140.4, 132.49, 133.32, 907849.86),
    @(1171, 45534.66666666666, 133.33, 134.95, 132.04, 134.47, 690059.313),
    @(1172, 45534.75, 134.47, 137.91, 134.4, 135.98, 441372.594),
    @(1173, 45534.83333333334, 135.98, 137.04, 135.65, 136.45, 120257.132),
    @(1174, 45534.91666666666, 136.45, 138.48, 136.42, 137.99, 200335.916),
    @(1175, 45535, 137.99, 139.23, 137.86, 137.95, 183688.248),
    @(1176, 45535.08333333334, 137.95, 138.81, 137.5, 138.29, 81273.049),
    @(1177, 45535.16666666666, 138.29, 138.77, 137.65, 137.66, 60952.884),
    @(1178, 45535.25, 137.66, 138.34, 137.38, 137.59, 69374.931),
    @(1179, 45535.33333333334, 137.59, 137.85, 135.79, 136.63, 178517.16),
    @(1180, 45535.41666666666, 136.63, 137.5, 136.49, 137.01, 121671.29),
    @(1181, 45535.5, 137.02, 137.29, 135.67, 135.82, 103703.121),
    @(1182, 45535.58333333334, 135.82, 136.66, 135.57, 136.63, 73326.893),
    @(1183, 45535.66666666666, 136.63, 136.7, 134.89, 135.73, 78017.689),
    @(1184, 45535.75, 135.73, 136.54, 134.49, 134.97, 100597.125),
    @(1185, 45535.83333333334, 134.97, 135.55, 134.51, 135.47, 86965.66099999999),
    @(1186, 45535.91666666666, 135.93, 136.22, 134.86, 135.35, 76562.641),
    @(1187, 45536, 135.35, 135.62, 133.5, 133.63, 96352.18799999999),
    @(1188, 45536.08333333334, 133.62, 134.46, 132.88, 134.07, 182954.106),
    @(1189, 45536.16666666666, 134.07, 134.09, 131.09, 131.59, 223355.854),
    @(1190, 45536.25, 131.59, 133, 131.49, 132.81, 173725.93),
    @(1191, 45536.33333333334, 132.81, 132.84, 130.78, 131.22, 235253.981),
    @(1192, 45536.41666666666, 131.2, 131.98, 130.18, 131.91, 201590.493),
    @(1193, 45536.5, 131.91, 132.16, 130.21, 131.58, 148239.468),
    @(1194, 45536.58333333334, 131.58, 134.87, 130, 134.29, 381278.873),
    @(1195, 45536.66666666666, 134.29, 134.75, 131.87, 132.65, 175884.671),
    @(1196, 45536.75, 132.65, 135.92, 132.44, 135.23, 187688.457),
    @(1197, 45536.83333333334, 135.24, 136.15, 133.12, 133.85, 282591.373),
    @(1198, 45536.91666666666, 133.86, 133.99, 127.66, 128.5, 639529.347),
    @(1199, 45537, 128.5, 130.48, 128.27, 129.74, 322492.184),
    @(1200, 45537.08333333334, 129.73, 130.7, 128.5, 130.21, 244653.936),
    @(1201, 45537.16666666666, 130.21, 130.37, 127.14, 129.09, 413646.424),
    @(1202, 45537.25, 129.08, 129.75, 126.8, 127.66, 226326.184),
    @(1203, 45537.33333333334, 127.65, 131.83, 127.22, 131.68, 454771.056),
    @(1204, 45537.41666666666, 131.68, 132.93, 131.27, 131.89, 330531.241),
    @(1205, 45537.5, 131.83, 132.47, 130.67, 130.77, 203499.908),
    @(1206, 45537.58333333334, 130.77, 133.45, 130, 132.54, 308001.345),
    @(1207, 45537.66666666666, 132.54, 133.66, 131.71, 133.34, 156958.792),
    @(1208, 45537.75, 133.34, 133.87, 132.51, 133.05, 129419.618),
    @(1209, 45537.83333333334, 133.04, 135.05, 132.81, 134.17, 191542.078),
    @(1210, 45537.91666666666, 134.17, 135.39, 133.99, 134.94, 159791.224),
    @(1211, 45538, 134.95, 135.47, 134.21, 135.44, 117606.885),
    @(1212, 45538.08333333334, 135.45, 138.13, 135.37, 135.51, 284290.817),
    @(1213, 45538.16666666666, 135.5, 135.89, 134.5, 135.02, 159864.489),
    @(1214, 45538.25, 135.02, 135.02, 133.54, 134.6, 177450.102),
    @(1215, 45538.33333333334, 134.6, 134.64, 132.77, 133.57, 169211.726),
    @(1216, 45538.41666666666, 133.58, 134.53, 133.44, 134.36, 110048.202),
    @(1217, 45538.5, 134.36, 135.5, 131.15, 131.58, 241880.939),
    @(1218, 45538.58333333334, 131.57, 131.57, 129.31, 129.8, 425322.328),
    @(1219, 45538.66666666666, 129.81, 131.5, 129.31, 130.16, 203036.496),
    @(1220, 45538.75, 130.17, 132.21, 129.74, 130.95, 171754.546),
    @(1221, 45538.83333333334, 130.96, 131.61, 129.31, 130.13, 211832.417),
    @(1222, 45538.91666666666, 130.13, 130.38, 127.15, 127.54, 355082.978),
    @(1223, 45539, 127.54, 128.89, 122.6, 126.83, 977639.1310000001),
    @(1224, 45539.08333333334, 126.84, 128.02, 126.21, 127.7, 222767.287),
    @(1225, 45539.16666666666, 127.7, 127.95, 126.55, 126.86, 157615.976),
    @(1226, 45539.25, 126.85, 131.44, 126.85, 130.3, 347500.821),
    @(1227, 45539.33333333334, 130.29, 130.71, 129.53, 129.59, 186282.719),
    @(1228, 45539.41666666666, 129.6, 129.85, 128.11, 128.43, 212624.822),
    @(1229, 45539.5, 128.43, 129.65, 127.88, 129.31, 265934.981),
    @(1230, 45539.58333333334, 129.31, 134.13, 128.67, 133.83, 766214.4570000001),
    @(1231, 45539.66666666666, 133.82, 134.98, 132.58, 132.61, 334861.194),
    @(1232, 45539.75, 132.62, 133.87, 132.58, 132.72, 124255.408)
)

$lastExistingRow = 1168
$srcDateCell = $ws.Cells.Item($lastExistingRow, 1)

foreach ($row in $newRows) {
    $r = $row[0]

    $dateCell = $ws.Cells.Item($r, 1)
    $srcDateCell.Copy($dateCell)
    $dateCell.Value = $row[1]

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
